$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old schedule info (group/time slots) that used to live in
# E2:G2, G3, E4:G4 and G5, keeping the cell styles intact.
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""

$ws.Range("G3").Value = ""

$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""

$ws.Range("G5").Value = ""

# New marker cell for row 4.
$ws.Range("B4").Value = "~"

# Update the saved selection.
$ws.Range("B7").Select()
